# Verizon location-search test: add a new "StoreLocator" sheet holding the
# zip code used to drive the store-locator search, and leave it as the
# active/selected sheet (with G1 highlighted, mirroring the saved UI state).

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab strip (becomes sheetId 3 / the 3rd tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "StoreLocator"

# Zip code value, stored as quote-prefixed text (so "11209" round-trips as
# a string rather than a number).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "'11209"

# Leave the new sheet active, with G1 selected.
$ws.Activate()
$ws.Range("G1").Select()
